$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "SELECT DISTINCT`n     std.dbgap_accession AS `"dbGaP Accession`",`n     std.study_name AS `"Study Name`"`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`nWHERE `n    std.dbgap_accession = 'phs001437' AND srv.first_event IN ('Not Reported')"
$ws.Range("C2").Value = "SELECT `n    COUNT(DISTINCT dgn.diagnosis) AS Diagnoses,`n    COUNT(DISTINCT prt.participant_id) AS Participants,`n    COUNT(DISTINCT std.study_id) AS Studies`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`nWHERE `n    std.dbgap_accession = 'phs001437' AND srv.first_event IN ('Not Reported');"
$ws.Range("B3").Value = "SELECT DISTINCT`n    prt.participant_id AS `"Participant Id`",`n    prt.race AS `"Race`",`n    prt.sex_at_birth AS `"Sex at Birth`",`n    std.dbgap_accession AS `"dbGaP Accession`"`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`nWHERE `n    std.dbgap_accession = 'phs001437' AND srv.first_event IN ('Not Reported')`nORDER BY `n    prt.participant_id ASC`nLIMIT 100;"
$ws.Range("B4").Value = "SELECT DISTINCT`n    prt.participant_id AS `"Participant Id`",`n    dgn.diagnosis_id AS `"Diagnosis Id`",`n    dgn.diagnosis AS `"Diagnosis`",`n    dgn.diagnosis_classification_system AS `"Diagnosis Classification System`",`n    dgn.diagnosis_basis AS `"Diagnosis Basis`",`n    dgn.tumor_classification AS `"Tumor Classification`",`n    dgn.anatomic_site AS `"Anatomic Site`",`n    CASE `n    WHEN dgn.age_at_diagnosis = -999 THEN 'Not Reported'`n    WHEN dgn.age_at_diagnosis >= 1000 THEN `n        substr(dgn.age_at_diagnosis, 1, length(dgn.age_at_diagnosis) - 3) || ',' || substr(dgn.age_at_diagnosis, -3)`n    ELSE `n        dgn.age_at_diagnosis `nEND AS `"Age at Diagnosis (days)`",`n    std.dbgap_accession AS `"dbGaP Accession`"`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`nWHERE `n    std.dbgap_accession = 'phs001437' AND srv.first_event IN ('Not Reported')`nORDER BY `n    prt.participant_id ASC`nLIMIT 100;"
$ws.Range("B5").Value = "SELECT`n    DISTINCT prt.participant_id AS `"Participant Id`",`n    trt.treatment_id AS `"Treatment Id`",`n    CASE `n    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'`n    WHEN trt.age_at_treatment_start >= 1000 THEN `n        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)`n    ELSE `n        trt.age_at_treatment_start `nEND AS `"Age at Treatment Start`",`n    CASE `n    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'`n    WHEN trt.age_at_treatment_end >= 1000 THEN `n        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)`n    ELSE `n        trt.age_at_treatment_end `nEND AS `"Age at Treatment End`",`n    trt.treatment_type AS `"Treatment Type`",`n    REPLACE(trt.treatment_agent, ';', ', ') AS `"Treatment Agent`",`n    std.dbgap_accession AS `"dbGaP Accession`"`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`nWHERE `n    std.dbgap_accession = 'phs001437' AND srv.first_event IN ('Not Reported')`nORDER BY `n    trt.treatment_id ASC`nLIMIT 100;"
$ws.Range("B6").Value = "SELECT DISTINCT`n    prt.participant_id AS `"Participant Id`",`n    trr.treatment_response_id AS `"Treatment Response Id`",`n    trr.response AS `"Response`",`n    CASE `n        WHEN trr.age_at_response = -999 THEN 'Not Reported'`n        WHEN trr.age_at_response >= 1000 THEN `n            substr(trr.age_at_response, 1, length(trr.age_at_response) - 3) || ',' || substr(trr.age_at_response, -3)`n        ELSE `n            trr.age_at_response `n    END AS `"Age at Response`",`n    trr.response_category AS `"Response Category`",`n    trr.response_system AS `"Response System`",`n    std.dbgap_accession AS `"dbGaP Accession`"`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`nWHERE `n   std.dbgap_accession = 'phs001437' AND srv.first_event IN ('Not Reported')`nORDER BY `n    prt.participant_id ASC`nLIMIT 100;"
$ws.Range("B7").Value = "SELECT DISTINCT`n    prt.participant_id AS `"Participant Id`",`n    srv.survival_id AS `"Survival Id`",`n    srv.last_known_survival_status AS `"Last Known Survival Status`",`n    CASE `n    WHEN srv.age_at_last_known_survival_status = -999 THEN 'Not Reported'`n    WHEN srv.age_at_last_known_survival_status >= 1000 THEN `n        substr(srv.age_at_last_known_survival_status, 1, length(srv.age_at_last_known_survival_status) - 3) || ',' || substr(srv.age_at_last_known_survival_status, -3)`n    ELSE `n        srv.age_at_last_known_survival_status `nEND AS `"Age at Last Known Survival Status`",`n    srv.first_event AS `"First Event`",`n    srv.cause_of_death AS `"Cause of Death`",`n    std.dbgap_accession AS `"dbGaP Accession`"`nFROM `n    df_study std`nLEFT JOIN `n    df_participant prt ON std.study_id = prt.`"study.study_id`"`nLEFT JOIN `n    df_diagnoses dgn ON prt.participant_id = dgn.`"participant.participant_id`"`nLEFT JOIN `n    df_treatments trt ON prt.participant_id = trt.`"participant.participant_id`"`nLEFT JOIN `n    df_treatment_resp trr ON prt.participant_id = trr.`"participant.participant_id`"`nLEFT JOIN `n    df_survival srv ON prt.participant_id = srv.`"participant.participant_id`"`nLEFT JOIN `n    df_reference_files rfs ON std.study_id = rfs.`"study.study_id`"`nWHERE `n    std.dbgap_accession = 'phs001437' AND srv.first_event IN ('Not Reported')`nORDER BY `n    prt.participant_id ASC`nLIMIT 100;"
